# Generate Report for Handback
# Update the timestamp text values recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file
$wsOverview.Range("G2").Value = "2016-08-21 19:07:57"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsZhCn.Range("H2").Value = "2016-08-21 19:07:53"
$wsZhCn.Range("K2").Value = "2016-08-21 19:08:12"

# de-de sheet: Correspond Handback DateTime for the first file
$wsDeDe.Range("K2").Value = "2016-08-21 19:08:19"
